$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 374
$ws.Range("I12").Value = 120
$ws.Range("K12").Value = 120
$ws.Range("M12").Value = 50

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1484.6
$ws.Range("I63").Value = 1484.6
$ws.Range("K63").Value = 1484.6
$ws.Range("M63").Value = -798.5999999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1484.6
$ws.Range("I66").Value = 1484.6
$ws.Range("K66").Value = 7423
$ws.Range("M66").Value = -3991

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4000
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2302.7
$ws.Range("I132").Value = 2280.7778
$ws.Range("K132").Value = 6842.3334
$ws.Range("M132").Value = -4312.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1944.9333
$ws.Range("I94").Value = 2089.5386
$ws.Range("K94").Value = 2089.5386
$ws.Range("M94").Value = -1638.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("K105").Value = 2000
$ws.Range("M105").Value = -253

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 857.1429000000001
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1000
$ws.Range("J4").Value = 1000
$ws.Range("L4").Value = 1000
$ws.Range("N4").Value = -1224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 536
$ws.Range("I16").Value = 475.33334
$ws.Range("J16").Value = 627
$ws.Range("K16").Value = 475.33334
$ws.Range("L16").Value = 627
$ws.Range("M16").Value = -188.33334
$ws.Range("N16").Value = -1201

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 4496.091
$ws.Range("I19").Value = 1861.6
$ws.Range("J19").Value = 6691.5
$ws.Range("K19").Value = 1861.6
$ws.Range("L19").Value = 6691.5
$ws.Range("M19").Value = -1691.6
$ws.Range("N19").Value = -7031.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 475
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 4496.091
$ws.Range("I24").Value = 1861.6
$ws.Range("J24").Value = 6691.5
$ws.Range("K24").Value = 1861.6
$ws.Range("L24").Value = 6691.5
$ws.Range("M24").Value = -1691.6
$ws.Range("N24").Value = -7031.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2896.7778
$ws.Range("I33").Value = 1178.5
$ws.Range("J33").Value = 6333.3335
$ws.Range("K33").Value = 1178.5
$ws.Range("L33").Value = 6333.3335
$ws.Range("M33").Value = -799.5
$ws.Range("N33").Value = -7091.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H48").Value = 152250
$ws.Range("J48").Value = 152250
$ws.Range("L48").Value = 152250
$ws.Range("N48").Value = -153202

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H56").Value = 4999.6665
$ws.Range("I56").Value = 4999.6665
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 4999.6665
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -4154.6665
$ws.Range("N56").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 549.9048
$ws.Range("I107").Value = 501.94446
$ws.Range("J107").Value = 837.6667
$ws.Range("K107").Value = 501.94446
$ws.Range("L107").Value = 837.6667
$ws.Range("M107").Value = 1418.05554
$ws.Range("N107").Value = -4677.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 536
$ws.Range("I113").Value = 475.33334
$ws.Range("J113").Value = 627
$ws.Range("K113").Value = 475.33334
$ws.Range("L113").Value = 627
$ws.Range("M113").Value = 1694.66666
$ws.Range("N113").Value = -4967

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 242.33333
$ws.Range("I29").Value = 118
$ws.Range("J29").Value = 366.66666
$ws.Range("K29").Value = 354
$ws.Range("L29").Value = 1099.99998
$ws.Range("M29").Value = -77
$ws.Range("N29").Value = -1653.99998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1175.6666
$ws.Range("I34").Value = 301.83334
$ws.Range("J34").Value = 2049.5
$ws.Range("K34").Value = 905.5000200000001
$ws.Range("L34").Value = 6148.5
$ws.Range("M34").Value = -821.5000200000001
$ws.Range("N34").Value = -6316.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 230.66667
$ws.Range("I55").Value = 230.66667
$ws.Range("K55").Value = 692.00001
$ws.Range("M55").Value = -515.00001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 20000
$ws.Range("J51").Value = 20000
$ws.Range("L51").Value = 20000
$ws.Range("N51").Value = -21018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2492.2354
$ws.Range("I102").Value = 2336.125
$ws.Range("K102").Value = 2336.125
$ws.Range("M102").Value = -714.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 1000
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1678.7142
$ws.Range("I113").Value = 1150.4
$ws.Range("J113").Value = 2999.5
$ws.Range("K113").Value = 1150.4
$ws.Range("L113").Value = 2999.5
$ws.Range("M113").Value = 1019.6
$ws.Range("N113").Value = -7339.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 15261.25
$ws.Range("I122").Value = 18682.334
$ws.Range("K122").Value = 56047.00199999999
$ws.Range("M122").Value = -53597.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5704.722
$ws.Range("I126").Value = 3257.0833
$ws.Range("J126").Value = 10600
$ws.Range("K126").Value = 9771.249899999999
$ws.Range("L126").Value = 31800
$ws.Range("M126").Value = -7301.249899999999
$ws.Range("N126").Value = -36740

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 727.95654
$ws.Range("I22").Value = 663.9524
$ws.Range("K22").Value = 663.9524
$ws.Range("M22").Value = -368.9524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 727.95654
$ws.Range("I27").Value = 663.9524
$ws.Range("K27").Value = 663.9524
$ws.Range("M27").Value = -556.9524

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2563.8572
$ws.Range("J46").Value = 1999
$ws.Range("L46").Value = 1999
$ws.Range("N46").Value = -2375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 326
$ws.Range("I113").Value = 326
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 978
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1192
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 79999
$ws.Range("J128").Value = 79999
$ws.Range("L128").Value = 79999
$ws.Range("N128").Value = -89959
